# Update "想去人数" (number of people interested) figures in the
# "展览" (Exhibition) and "全部类型" (All types) sheets.
# Sheets "演出" and "本地生活" only contain header rows, so they are untouched.

$wb = $excel.ActiveWorkbook

$targetSheetNames = @("展览", "全部类型")

# New values for column F, rows 2-11 are identical across both sheets.
$commonValues = @{
    2  = 1503
    3  = 23
    4  = 968
    5  = 62
    6  = 2315
    7  = 38
    8  = 1419
    9  = 66
    10 = 154
    11 = 49
}

foreach ($sheetName in $targetSheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)

    foreach ($row in $commonValues.Keys) {
        $ws.Cells.Item($row, 6).Value = $commonValues[$row]
    }

    # Row 12 ends at the same value (371) on both sheets even though the
    # starting values differed (349 on 展览, 350 on 全部类型).
    $ws.Cells.Item(12, 6).Value = 371
}
